$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Email" column only held real data in rows 1 (header) and 2 (the sample
# record). Rather than deleting the whole column B (which would shift every
# row and disturb the sparse data used throughout the rest of the sheet),
# shift just the cells in rows 1 and 2 from column C onward one column to the
# left (into B..I), then clear out the now-unused column J, for those two
# rows only.
foreach ($r in 1,2) {
    for ($c = 2; $c -le 9; $c++) {
        $srcVal = $ws.Cells.Item($r, $c + 1).Value()
        $ws.Cells.Item($r, $c).Value = $srcVal
    }
    $ws.Cells.Item($r, 10).ClearContents()
}

# The shifted-in header/email cell no longer needs the Hyperlink style.
$ws.Cells.Item(2, 2).Style = "Normal"

# Remove the now-orphaned mailto hyperlink itself.
if ($ws.Hyperlinks.Count -gt 0) {
    $ws.Hyperlinks.Delete()
}

# The built-in "Hyperlink" cell style is no longer used anywhere in the
# workbook now that the emailed-styled cell is gone - drop it so the style
# table collapses back down to just "Normal".
$hyperlinkStyle = $wb.Styles.Item("Hyperlink")
$hyperlinkStyle.Delete()

# Update selection to match the target workbook (active cell D5).
$ws.Range("D5").Select()

$wb.Save()
